$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 3
$ws.Range("H3").Value = 44938.4
$ws.Range("J3").Value = 44938.4
$ws.Range("L3").Value = 44938.4
$ws.Range("N3").Value = -45166.4
# row 32
$ws.Range("H32").Value = 1055
$ws.Range("I32").Value = 999
$ws.Range("K32").Value = 999
$ws.Range("M32").Value = -673
# row 88
$ws.Range("H88").Value = 2643
$ws.Range("I88").Value = 3654
$ws.Range("K88").Value = 3654
$ws.Range("M88").Value = -3248
# row 91
$ws.Range("H91").Value = 2643
$ws.Range("I91").Value = 3654
$ws.Range("K91").Value = 3654
$ws.Range("M91").Value = -2250
# row 102
$ws.Range("H102").Value = 44938.4
$ws.Range("J102").Value = 44938.4
$ws.Range("L102").Value = 44938.4
$ws.Range("N102").Value = -51428.4
# row 103
$ws.Range("H103").Value = 5945.4165
$ws.Range("J103").Value = 6408.636
$ws.Range("L103").Value = 19225.908
$ws.Range("N103").Value = -20397.908
# row 127
$ws.Range("H127").Value = 2885.8
$ws.Range("I127").Value = 2857.5
$ws.Range("J127").Value = 2999
$ws.Range("K127").Value = 8572.5
$ws.Range("L127").Value = 8997
$ws.Range("M127").Value = -3612.5
$ws.Range("N127").Value = -18917
# row 132
$ws.Range("H132").Value = 31271.555
$ws.Range("I132").Value = 42810.8
$ws.Range("K132").Value = 128432.4
$ws.Range("M132").Value = -125902.4
# row 138
$ws.Range("H138").Value = 2362.75
$ws.Range("I138").Value = 1186
$ws.Range("J138").Value = 10600
$ws.Range("K138").Value = 3558
$ws.Range("L138").Value = 31800
$ws.Range("M138").Value = 1582
$ws.Range("N138").Value = -42080

$ws = $wb.Worksheets.Item("ARM")
# row 28
$ws.Range("H28").Value = 50000
$ws.Range("I28").Value = 50000
$ws.Range("K28").Value = 50000
$ws.Range("M28").Value = -49808
# row 32
$ws.Range("H32").Value = 2943558.5
$ws.Range("I32").Value = 551.4483
$ws.Range("K32").Value = 551.4483
$ws.Range("M32").Value = -264.4483
# row 45
$ws.Range("H45").Value = 2699.889
$ws.Range("I45").Value = 2036.25
$ws.Range("K45").Value = 2036.25
$ws.Range("M45").Value = -1659.25
# row 61
$ws.Range("H61").Value = 2126.3076
$ws.Range("I61").Value = 2012.909
$ws.Range("K61").Value = 2012.909
$ws.Range("M61").Value = -1800.909
# row 99
$ws.Range("H99").Value = 50000
$ws.Range("I99").Value = 50000
$ws.Range("K99").Value = 50000
$ws.Range("M99").Value = -47005
# row 102
$ws.Range("H102").Value = 4810143
$ws.Range("I102").Value = 6250686
$ws.Range("K102").Value = 6250686
$ws.Range("M102").Value = -6249064
# row 122
$ws.Range("H122").Value = 2186.5833
$ws.Range("I122").Value = 1577.5555
$ws.Range("K122").Value = 4732.666499999999
$ws.Range("M122").Value = -2282.666499999999
# row 136
$ws.Range("H136").Value = 2126.3076
$ws.Range("I136").Value = 2012.909
$ws.Range("K136").Value = 6038.727000000001
$ws.Range("M136").Value = -3488.727000000001

$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 71429720
$ws.Range("I99").Value = 90910250
$ws.Range("J99").Value = 1101.3334
$ws.Range("K99").Value = 90910250
$ws.Range("L99").Value = 1101.3334
$ws.Range("M99").Value = -90908752
$ws.Range("N99").Value = -4097.3334
# row 105
$ws.Range("H105").Value = 4785802
$ws.Range("I105").Value = 6993956.5
$ws.Range("K105").Value = 6993956.5
$ws.Range("M105").Value = -6992209.5
# row 106
$ws.Range("H106").Value = 32499.666
$ws.Range("J106").Value = 32499.666
$ws.Range("L106").Value = 32499.666
$ws.Range("N106").Value = -35023.666
# row 107
$ws.Range("H107").Value = 25005070
$ws.Range("I107").Value = 55557490
$ws.Range("J107").Value = 7636.1816
$ws.Range("K107").Value = 55557490
$ws.Range("L107").Value = 7636.1816
$ws.Range("M107").Value = -55555570
$ws.Range("N107").Value = -11476.1816

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1774
# row 22
$ws.Range("H22").Value = 1258.125
$ws.Range("I22").Value = 748.5
$ws.Range("J22").Value = 2107.5
$ws.Range("K22").Value = 748.5
$ws.Range("L22").Value = 2107.5
$ws.Range("M22").Value = -398.5
$ws.Range("N22").Value = -2807.5
# row 58
$ws.Range("H58").Value = 2927.652
$ws.Range("I58").Value = 1702.9474
$ws.Range("K58").Value = 1702.9474
$ws.Range("M58").Value = -1499.9474
# row 86
$ws.Range("H86").Value = 5999.5
$ws.Range("I86").Value = 5999.5
$ws.Range("K86").Value = 5999.5
$ws.Range("M86").Value = -4876.5
# row 89
$ws.Range("H89").Value = 5999.5
$ws.Range("I89").Value = 5999.5
$ws.Range("K89").Value = 29997.5
$ws.Range("M89").Value = -24381.5
# row 99
$ws.Range("H99").Value = 2908.1667
$ws.Range("I99").Value = 2899.818
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2899.818
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1401.818
$ws.Range("N99").Value = -5996
# row 105
$ws.Range("H105").Value = 871.1667
$ws.Range("I105").Value = 647.8
$ws.Range("J105").Value = 1988
$ws.Range("K105").Value = 647.8
$ws.Range("L105").Value = 1988
$ws.Range("M105").Value = 1099.2
$ws.Range("N105").Value = -5482
# row 107
$ws.Range("H107").Value = 863.44446
$ws.Range("I107").Value = 258.66666
$ws.Range("K107").Value = 258.66666
$ws.Range("M107").Value = 1661.33334
# row 113
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5540
# row 126
$ws.Range("H126").Value = 2908.1667
$ws.Range("I126").Value = 2899.818
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8699.454000000002
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6229.454000000002
$ws.Range("N126").Value = -13940
# row 132
$ws.Range("H132").Value = 2011.6666
$ws.Range("I132").Value = 1983.3334
$ws.Range("J132").Value = 2266.6667
$ws.Range("K132").Value = 5950.0002
$ws.Range("L132").Value = 6800.000100000001
$ws.Range("M132").Value = -3420.0002
$ws.Range("N132").Value = -11860.0001
# row 136
$ws.Range("H136").Value = 2927.652
$ws.Range("I136").Value = 1702.9474
$ws.Range("K136").Value = 5108.8422
$ws.Range("M136").Value = -2558.8422

$ws = $wb.Worksheets.Item("CUL")
# row 26
$ws.Range("H26").Value = 317.125
$ws.Range("I26").Value = 276.7143
$ws.Range("J26").Value = 600
$ws.Range("K26").Value = 830.1428999999999
$ws.Range("L26").Value = 1800
$ws.Range("M26").Value = -542.1428999999999
$ws.Range("N26").Value = -2376

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("I80").Value = 9999
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 9999
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -9001
$ws.Range("N80").ClearContents()
# row 83
$ws.Range("I83").Value = 9999
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 49995
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -45003
$ws.Range("N83").ClearContents()
# row 102
$ws.Range("H102").Value = 2160.647
$ws.Range("I102").Value = 2160.647
$ws.Range("K102").Value = 2160.647
$ws.Range("M102").Value = -538.6469999999999
# row 122
$ws.Range("H122").Value = 2612.6667
$ws.Range("I122").Value = 1298
$ws.Range("K122").Value = 3894
$ws.Range("M122").Value = -1444
# row 132
$ws.Range("H132").Value = 26777.904
$ws.Range("I132").Value = 31302.371
$ws.Range("K132").Value = 93907.113
$ws.Range("M132").Value = -91377.113
# row 140
$ws.Range("H140").Value = 131110.8
$ws.Range("J140").Value = 93964.25
$ws.Range("L140").Value = 93964.25
$ws.Range("N140").Value = -104324.25

$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3471.8572
$ws.Range("I7").Value = 3004
$ws.Range("J7").Value = 3549.8333
$ws.Range("K7").Value = 3004
$ws.Range("L7").Value = 3549.8333
$ws.Range("M7").Value = -2892
$ws.Range("N7").Value = -3773.8333
# row 40
$ws.Range("H40").Value = 3478
$ws.Range("I40").Value = 3478
$ws.Range("K40").Value = 3478
$ws.Range("M40").Value = -3342
# row 122
$ws.Range("H122").Value = 3000
$ws.Range("J122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
# row 125
$ws.Range("H125").Value = 15000
$ws.Range("J125").Value = 15000
$ws.Range("L125").Value = 15000
$ws.Range("N125").Value = -24840
# row 126
$ws.Range("H126").Value = 3471.8572
$ws.Range("I126").Value = 3004
$ws.Range("J126").Value = 3549.8333
$ws.Range("K126").Value = 9012
$ws.Range("L126").Value = 10649.4999
$ws.Range("M126").Value = -6542
$ws.Range("N126").Value = -15589.4999
# row 136
$ws.Range("H136").Value = 2148.3333
$ws.Range("I136").Value = 1962.625
$ws.Range("K136").Value = 5887.875
$ws.Range("M136").Value = -3337.875

$ws = $wb.Worksheets.Item("WVR")
# row 100
$ws.Range("H100").Value = 1756.5714
$ws.Range("I100").Value = 1756.5714
$ws.Range("K100").Value = 3513.1428
$ws.Range("M100").Value = -2972.1428
# row 132
$ws.Range("H132").Value = 2244.6365
$ws.Range("I132").Value = 1965.7778
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 5897.3334
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -3367.3334
$ws.Range("N132").Value = -15558.5

Write-Output "Applied market price updates across sheets."
